# Simulation_Output_Settings.xlsx update
#
# Adds two more selectable "Total Energy Output Variable" rows to the
# TotalEnergy sheet and one more meter row to the Meters sheet (matching
# the new entries already present in the Sources lookup lists), then
# leaves the selection/active sheet where the author left it.

$wb = $excel.ActiveWorkbook

# --- TotalEnergy sheet: two additional output-variable choices ---
$wsTotal = $wb.Worksheets.Item("TotalEnergy")
$wsTotal.Range("A4").Value = "Electricity Interior Lighting [GJ]"
$wsTotal.Range("A5").Value = "Electricity Fans [GJ]"
$wsTotal.Range("A9").Select()

# --- Meters sheet: one additional meter/frequency pair ---
$wsMeters = $wb.Worksheets.Item("Meters")
$wsMeters.Range("A4").Value = "Cooling:Electricity"
$wsMeters.Range("B4").Value = "Monthly"

# Meters stays the active/visible sheet, with the selection parked at B6
$wsMeters.Activate()
$wsMeters.Range("B6").Select()
